$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header row: "_old" suffix -> "_FV2304", "_new" suffix -> "_FV2310"
# (columns A:J use the "_old" variant, column K is "diff", columns L:U use the "_new" variant)
$headerRange = $ws.Range("A1:U1")
$headerRange.Replace("_old", "_FV2304", [Microsoft.Office.Interop.Excel.XlLookAt]::xlPart)
$headerRange.Replace("_new", "_FV2310", [Microsoft.Office.Interop.Excel.XlLookAt]::xlPart)

# Turn the used range into a real Excel Table (ListObject) with an AutoFilter
$dataRange = $ws.Range("A1:U83")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $dataRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

# Freeze the header row (split/freeze pane under row 1)
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
